$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'320.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-3.53%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'42.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-6.70%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.152"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-9.22%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08151"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'4.333"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-3.33%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.780"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-12.81%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9516"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-4.04%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1117"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-3.32%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1850"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-4.16%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09364"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-5.90%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04621"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.19%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'7.428"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-28.67%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'0.22%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001286"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.29%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006025"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.42%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-0.19%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E19").Value = "'-0.03%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1389"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.83%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2630"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.87%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-0.40%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001252"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-4.40%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004317"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-6.78%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001113"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-13.22%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0002988"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-20.30%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02596"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-7.33%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05535"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-4.28%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007835"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.22%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1393"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.93%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.006620"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-9.13%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002120"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.38%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008459"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-6.54%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3461"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.59%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007005"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-5.09%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.06%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003474"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-0.92%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003541"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.98%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.06%"
$ws.Range("E51").Style = "Normal"
